$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9: s2s3 (added first so its string gets the next shared-string index)
$ws.Range("A9").Value = "s2s3"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 1

# Add the new "cosimbug" marker to row 8 (nufft512), column F
$ws.Range("F8").Value = "cosimbug"

# Add new row 10: merge7
$ws.Range("A10").Value = "merge7"
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 16
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 7

# Update selection to reflect the new active cell after entry (B11)
$ws.Range("B11").Select()
